$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 106, shifting existing rows 106-156 down to 107-157.
$ws.Rows("106:106").Insert()

# Populate the newly inserted row 106 with the new record's data.
$ws.Range("A106").Value = 3
$ws.Range("B106").Value = "Femacal de La Calera"
$ws.Range("C106").Value = "Coquimbo"
$ws.Range("D106").Value = 44529
$ws.Range("E106").Value = 5
$ws.Range("F106").Value = 100112010
$ws.Range("G106").Value = "Achicoria"
$ws.Range("H106").Value = "Sin especificar"
$ws.Range("I106").Value = "Primera"
$ws.Range("J106").Value = 50
$ws.Range("K106").Value = 5500
$ws.Range("L106").Value = 5500
$ws.Range("M106").Value = 5500
$ws.Range("N106").Value = "$/caja 16 unidades"
$ws.Range("O106").Value = "Provincia de Quillota"
$ws.Range("P106").Value = 344
$ws.Range("Q106").Value = 16
$ws.Range("R106").Value = "Hortaliza"

# Make sure the date column keeps the date number format (style index 2).
$ws.Range("D106").NumberFormat = $ws.Range("D107").NumberFormat
